$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B11").Value = -0
$ws.Range("B21").Value = -0
$ws.Range("B31").Value = -5000000
$ws.Range("B34").Value = -3882857.142857143
$ws.Range("B35").Value = -3882857.142857143
